$wb = $excel.ActiveWorkbook

# Avoid confirmation prompts (e.g. when deleting a worksheet with data)
$excel.DisplayAlerts = $false

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wsDesarquivamentos = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$wsDesarquivamentos.Delete() | Out-Null

$excel.DisplayAlerts = $true
